$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.297.88"

$ws.Range("D3").Value = "3.554.93"
$ws.Range("E3").Value = "  +6.15%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'190.28"
$ws.Range("E5").Value = "  +9.22%  "

$ws.Range("D6").Value = "'564.49"
$ws.Range("E6").Value = "  +6.82%  "

$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  +4.58%  "

$ws.Range("D8").Value = "3.544.43"
$ws.Range("E8").Value = "  +5.67%  "

$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "'0.637"
$ws.Range("E10").Value = "  +4.49%  "

$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  +14.17%  "

$ws.Range("D12").Value = "'55.11"
$ws.Range("E12").Value = "  +2.54%  "

$ws.Range("D13").Value = "'0.0000274"
$ws.Range("E13").Value = "  +6.68%  "

$ws.Range("D14").Value = "'9.38"
$ws.Range("E14").Value = "  +2.89%  "

$ws.Range("D15").Value = "4.096.80"
$ws.Range("E15").Value = "  +5.77%  "

$ws.Range("D16").Value = "3.550.39"
$ws.Range("E16").Value = "  +6.34%  "

$ws.Range("D17").Value = "'18.66"
$ws.Range("E17").Value = "  +6.10%  "

$ws.Range("E18").Value = "  +3.56%  "

$ws.Range("D19").Value = "67.258.32"
$ws.Range("E19").Value = "  +5.18%  "

$ws.Range("D20").Value = "'12.12"
$ws.Range("E20").Value = "  +7.80%  "

$ws.Range("E21").Value = "  +3.81%  "

$ws.Range("D22").Value = "'423.40"
$ws.Range("E22").Value = "  +13.05%  "

$ws.Range("D23").Value = "'4.14"
$ws.Range("E23").Value = "  +11.11%  "

$ws.Range("D24").Value = "'85.47"
$ws.Range("E24").Value = "  +4.69%  "

$ws.Range("E25").Value = "  +1.95%  "

$ws.Range("D26").Value = "'11.10"
$ws.Range("E26").Value = "  -4.43%  "

$ws.Range("E27").Value = "  +8.72%  "

$ws.Range("D28").Value = "'12.35"
$ws.Range("E28").Value = "  +8.99%  "

$ws.Range("D29").Value = "'6.07"
$ws.Range("E29").Value = "  +0.63%  "

$ws.Range("D30").Value = "'9.12"
$ws.Range("E30").Value = "  +10.12%  "

$ws.Range("D31").Value = "'30.57"
$ws.Range("E31").Value = "  +5.65%  "

$ws.Range("D32").Value = "'633.26"
$ws.Range("E32").Value = "  -0.28%  "

$ws.Range("D33").Value = "'6.70"
$ws.Range("E33").Value = "  +4.05%  "

$ws.Range("D34").Value = "'11.79"
$ws.Range("E34").Value = "  +4.89%  "

$ws.Range("E35").Value = "  +5.58%  "

$ws.Range("D36").Value = "'60.43"
$ws.Range("E36").Value = "  +4.21%  "

$ws.Range("D37").Value = "0.0₃0832"
$ws.Range("E37").Value = "  +14.52%  "

$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'38.60"
$ws.Range("E38").Value = "  +5.63%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.149"
$ws.Range("E39").Value = "  +19.85%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").Value = "'0.390"
$ws.Range("E41").Value = "  +2.42%  "

$ws.Range("D42").Value = "'3.37"
$ws.Range("E42").Value = "  +11.26%  "

$ws.Range("D43").Value = "3.145.63"
$ws.Range("E43").Value = "  +5.59%  "

$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").Value = "'2.66"
$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("D46").Value = "'2.89"
$ws.Range("E46").Value = "  +10.28%  "

$ws.Range("D47").Value = "'3.36"
$ws.Range("E47").Value = "  +11.34%  "

$ws.Range("D48").Value = "'0.0420"
$ws.Range("E48").Value = "  +5.91%  "

$ws.Range("E49").Value = "  +2.11%  "

$ws.Range("E50").Value = "  +5.52%  "

$ws.Range("D51").Value = "'8.66"
$ws.Range("E51").Value = "  +8.96%  "

